$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''66.186.76'
$ws.Range('E2').Value = '''  -0.37%  '
$ws.Range('D3').Value = '''3.282.16'
$ws.Range('E3').Value = '''  -0.94%  '
$ws.Range('E4').Value = '''  +0.03%  '
$ws.Range('D5').Value = '''586.39'
$ws.Range('E5').Value = '''  +2.58%  '
$ws.Range('D6').Value = '''181.19'
$ws.Range('E6').Value = '''  +0.13%  '
$ws.Range('D7').Value = '''0.647'
$ws.Range('E7').Value = '''  +8.12%  '
$ws.Range('E8').Value = '''  +0.04%  '
$ws.Range('E9').Value = '''  -2.62%  '
$ws.Range('D10').Value = '''6.74'
$ws.Range('E10').Value = '''  +1.52%  '
$ws.Range('E11').Value = '''  +0.29%  '
$ws.Range('D12').Value = '''3.856.66'
$ws.Range('E12').Value = '''  -0.77%  '
$ws.Range('E13').Value = '''  -5.33%  '
$ws.Range('D14').Value = '''66.211.83'
$ws.Range('E14').Value = '''  -0.47%  '
$ws.Range('D15').Value = '''26.49'
$ws.Range('E15').Value = '''  -1.95%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '''0.0000163'
$ws.Range('E16').Value = '''  -1.85%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '''3.233.20'
$ws.Range('E17').Value = '''  -2.35%  '
$ws.Range('D18').Value = '''435.13'
$ws.Range('E18').Value = '''  -0.60%  '
$ws.Range('E19').Value = '''  -2.93%  '
$ws.Range('D20').Value = '''5.50'
$ws.Range('E20').Value = '''  -3.08%  '
$ws.Range('D21').Value = '''7.44'
$ws.Range('E21').Value = '''  -3.08%  '
$ws.Range('D22').Value = '''72.11'
$ws.Range('E22').Value = '''  -2.40%  '
$ws.Range('E23').Value = '''  +0.06%  '
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').Value = '''5.73'
$ws.Range('E24').Value = '''  +1.56%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '''3.422.92'
$ws.Range('E25').Value = '''  -1.06%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').Value = '''0.511'
$ws.Range('E26').Value = '''  -0.40%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = '''0.0000114'
$ws.Range('E27').Value = '''  -4.13%  '
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').Value = '''0.196'
$ws.Range('E28').Value = '''  +3.13%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '''8.87'
$ws.Range('E29').Value = '''  -1.06%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '''  +0.03%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''1.96'
$ws.Range('E31').Value = '''  +0.03%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '''22.32'
$ws.Range('E32').Value = '''  -2.16%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').Value = '''0.999'
$ws.Range('E33').Value = '''  +0.12%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '''5.20'
$ws.Range('E34').Value = '''  -1.84%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '''6.64'
$ws.Range('E35').Value = '''  -1.76%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').Value = '''1.19'
$ws.Range('E36').Value = '''  -1.86%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '''158.10'
$ws.Range('E37').Value = '''  -1.30%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '''1.42'
$ws.Range('E38').Value = '''  -5.24%  '
$ws.Range('D39').Value = '''26.33'
$ws.Range('E39').Value = '''  -3.54%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '''1.78'
$ws.Range('E40').Value = '''  -3.27%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '''2.794.03'
$ws.Range('E41').Value = '''  -0.41%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').Value = '''0.773'
$ws.Range('E42').Value = '''  -1.52%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '''4.35'
$ws.Range('E43').Value = '''  -2.51%  '
$ws.Range('D44').Value = '''40.26'
$ws.Range('E44').Value = '''  +0.22%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '''6.12'
$ws.Range('E45').Value = '''  -0.93%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '''0.0662'
$ws.Range('E46').Value = '''  -1.76%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = '''2.32'
$ws.Range('E47').Value = '''  -0.73%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '''319.53'
$ws.Range('E48').Value = '''  +0.30%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '''23.30'
$ws.Range('E49').Value = '''  -3.39%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '''0.0268'
$ws.Range('E50').Value = '''  -0.95%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '''0.105'
$ws.Range('E51').Value = '''  +5.51%  '
